$d = $word.ActiveDocument

# The title paragraph ("Student application") currently owns the hidden
# "_GoBack" bookmark at its end. In the target document that bookmark is
# relocated to the end of the new "professor" paragraph further down, so
# remove it here first (it will be re-inserted via the OOXML fragment below).
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
    # no-op if the bookmark isn't present
}

# Collapsed range positioned right after the existing title run, i.e. at the
# end of paragraph 1 (just before its paragraph mark).
$titlePara = $d.Paragraphs(1)
$insertionPoint = $d.Range($titlePara.Range.End - 1, $titlePara.Range.End - 1)

$newParagraphsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A student will have one or more subjects they attend. Each subject will be taught by a professor. A professor can have one or more subjects. Upon finishing a subject, the student will get a grade related to that subject. The grades will be given by the professor.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A student will have a name, surname, age, grade, department and other info fields. The student will also have a list of subjects they are attending and a list of grades that they have gotten.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A professor will have a name, surname, age, title, department and other info fields. A professor will have a list of subjects they teach.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A subject will have a title, a year at which it is taught, a description and other info fields.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A grade will have a grade number, a name of the subject it is related to, the name of the professor who gave it, the name of the student that got it and the overall final score.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($newParagraphsXml)
